# Apply the "promos" content refresh:
#  - Column B ("หัวcard") now uses a single unified promo name for every row.
#  - Column C ("ประเภทโปรโมชั่น") wording is shortened for the monthly/weekly rows.
#  - Column H ("โบนัส") bonus text is refreshed for the three AIS SUPER WiFi rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (รายวัน)
$ws.Range("B2").Value = "AIS เน็ตไม่จำกัด"

# Row 3 (รายเดือน)
$ws.Range("C3").Value = "เดือน"

# Row 4 (รายสัปดาห์)
$ws.Range("C4").Value = "สัปดาห์"

# Bonus text refresh
$ws.Range("H2").Value = "ใช้ครบ 5GB ลดเหลือ 64Kbps"
$ws.Range("H3").Value = "ใช้ครบ 5GB ลดเหลือ 64Kbps"
$ws.Range("H4").Value = "ใช้ครบ 5GB ลดเหลือ 64Kbps"

# Finish filling in column B for the remaining rows
$ws.Range("B3").Value = "AIS เน็ตไม่จำกัด"
$ws.Range("B4").Value = "AIS เน็ตไม่จำกัด"
$ws.Range("B5").Value = "AIS เน็ตไม่จำกัด"

# Row 5 (รายเดือน / โทรฟรี AIS)
$ws.Range("C5").Value = "รายเดือน"

# Column B got wider once every row shares the longer "AIS เน็ตไม่จำกัด" label.
$ws.Columns.Item(2).ColumnWidth = 12.86

# Reset the view's active cell/selection back to A1.
$ws.Range("A1").Select()
